$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.80606733943071
$ws.Range("C2").Value = 8.843668143971094
$ws.Range("E2").Value = 11.43129748758444
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 28.86646142076297
$ws.Range("H2").Value = 13.62671500601549
$ws.Range("M2").Value = 14.6231630645241
$ws.Range("N2").Value = 16.757642252731
$ws.Range("B3").Value = 12.2146141055841
$ws.Range("C3").Value = 8.372209569239471
$ws.Range("E3").Value = 11.32154172202978
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 28.44596566184433
$ws.Range("H3").Value = 13.6433237711046
$ws.Range("M3").Value = 14.34378902573054
$ws.Range("N3").Value = 16.83279564638506
$ws.Range("B4").Value = 11.839279386063
$ws.Range("C4").Value = 8.067154450047489
$ws.Range("E4").Value = 11.2575128469396
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 28.19935661521504
$ws.Range("H4").Value = 13.6574316886959
$ws.Range("M4").Value = 14.17367221208957
$ws.Range("N4").Value = 16.88094530113196
$ws.Range("B5").Value = 11.68349411096636
$ws.Range("C5").Value = 7.938981897473832
$ws.Range("E5").Value = 11.2322913619066
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 28.10191440196221
$ws.Range("H5").Value = 13.66415963763008
$ws.Range("M5").Value = 14.10480562971849
$ws.Range("N5").Value = 16.90107281190421
$ws.Range("B6").Value = 11.65746230198552
$ws.Range("C6").Value = 7.917467718001774
$ws.Range("E6").Value = 11.22815665215825
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 28.08592255289495
$ws.Range("H6").Value = 13.66533579237848
$ws.Range("M6").Value = 14.09340095375758
$ws.Range("N6").Value = 16.90444559540996
$ws.Range("B7").Value = 11.8371895564434
$ws.Range("C7").Value = 8.06544141501576
$ws.Range("E7").Value = 11.25716914318651
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 28.19802993357241
$ws.Range("H7").Value = 13.65751846699355
$ws.Range("M7").Value = 14.17274146835086
$ws.Range("N7").Value = 16.88121469595884
$ws.Range("B8").Value = 12.60479009118046
$ws.Range("C8").Value = 8.684372813020053
$ws.Range("E8").Value = 11.3927719149071
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 28.71915946231289
$ws.Range("H8").Value = 13.63162768980125
$ws.Range("M8").Value = 14.52660014518972
$ws.Range("N8").Value = 16.78314013219943
$ws.Range("B9").Value = 14.00488193862422
$ws.Range("C9").Value = 9.772572152309266
$ws.Range("E9").Value = 11.68414589094893
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 29.82631194644381
$ws.Range("H9").Value = 13.61205507976308
$ws.Range("M9").Value = 15.22742984469142
$ws.Range("N9").Value = 16.60663721249328
$ws.Range("B10").Value = 14.95992535911748
$ws.Range("C10").Value = 10.49364554644804
$ws.Range("E10").Value = 11.91199356223455
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 30.68237970619096
$ws.Range("H10").Value = 13.6169027298408
$ws.Range("M10").Value = 15.74086184431824
$ws.Range("N10").Value = 16.4864785514996
$ws.Range("B11").Value = 15.37689148558792
$ws.Range("C11").Value = 10.80439982165396
$ws.Range("E11").Value = 12.01824061885331
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 31.07904487514563
$ws.Range("H11").Value = 13.62331654315149
$ws.Range("M11").Value = 15.97297471153803
$ws.Range("N11").Value = 16.43385477296913
$ws.Range("B12").Value = 15.53217031265082
$ws.Range("C12").Value = 10.91957759190762
$ws.Range("E12").Value = 12.05881257419793
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 31.23012138008681
$ws.Range("H12").Value = 13.62635224381797
$ws.Range("M12").Value = 16.06057157754637
$ws.Range("N12").Value = 16.41421838739014
$ws.Range("B13").Value = 15.49884609178197
$ws.Range("C13").Value = 10.89488329926368
$ws.Range("E13").Value = 12.05006018642834
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 31.19754846905098
$ws.Range("H13").Value = 13.6256714366662
$ws.Range("M13").Value = 16.04172072874266
$ws.Range("N13").Value = 16.41843451664006
$ws.Range("B14").Value = 15.38971934737822
$ws.Range("C14").Value = 10.8139257040124
$ws.Range("E14").Value = 12.02157191664828
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 31.09145761404485
$ws.Range("H14").Value = 13.62355412047279
$ws.Range("M14").Value = 15.98018785871658
$ws.Range("N14").Value = 16.43223345209344
$ws.Range("B15").Value = 15.32253247881995
$ws.Range("C15").Value = 10.76401113224605
$ws.Range("E15").Value = 12.00416506192573
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 31.02658186017946
$ws.Range("H15").Value = 13.62233628466548
$ws.Range("M15").Value = 15.94245561600277
$ws.Range("N15").Value = 16.44072355402918
$ws.Range("B16").Value = 14.93231458764674
$ws.Range("C16").Value = 10.47298817498787
$ws.Range("E16").Value = 11.90509935823349
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 30.65658751909357
$ws.Range("H16").Value = 13.61656842180014
$ws.Range("M16").Value = 15.72565568390891
$ws.Range("N16").Value = 16.4899584730067
$ws.Range("B17").Value = 14.68837101025253
$ws.Range("C17").Value = 10.29002042543954
$ws.Range("E17").Value = 11.84496638107041
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 30.43133795815447
$ws.Range("H17").Value = 13.61410937045303
$ws.Range("M17").Value = 15.59221856977598
$ws.Range("N17").Value = 16.52068292755139
$ws.Range("B18").Value = 14.54642000919772
$ws.Range("C18").Value = 10.18315767823434
$ws.Range("E18").Value = 11.81062616323449
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 30.30247168264626
$ws.Range("H18").Value = 13.61309105522908
$ws.Range("M18").Value = 15.51533774480437
$ws.Range("N18").Value = 16.53854666161331
$ws.Range("B19").Value = 14.49807938457069
$ws.Range("C19").Value = 10.1466973182814
$ws.Range("E19").Value = 11.79904258332817
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 30.25896379368854
$ws.Range("H19").Value = 13.61281423321993
$ws.Range("M19").Value = 15.48928750337169
$ws.Range("N19").Value = 16.54462801936346
$ws.Range("B20").Value = 14.71450986554351
$ws.Range("C20").Value = 10.30966591398749
$ws.Range("E20").Value = 11.85134236124788
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 30.45524586573362
$ws.Range("H20").Value = 13.61433013501253
$ws.Range("M20").Value = 15.60643744222952
$ws.Range("N20").Value = 16.51739241547313
$ws.Range("B21").Value = 15.42184428007869
$ws.Range("C21").Value = 10.83777281031928
$ws.Range("E21").Value = 12.0299307072291
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 31.12259687596976
$ws.Range("H21").Value = 13.62415954477673
$ws.Range("M21").Value = 15.99827037541053
$ws.Range("N21").Value = 16.42817248443596
$ws.Range("B22").Value = 15.86883713935812
$ws.Range("C22").Value = 11.16835614071169
$ws.Range("E22").Value = 12.14860500509405
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 31.5637306337529
$ws.Range("H22").Value = 13.63412166309906
$ws.Range("M22").Value = 16.25257293601814
$ws.Range("N22").Value = 16.37155799960344
$ws.Range("B23").Value = 15.63169734883563
$ws.Range("C23").Value = 10.9932542205292
$ws.Range("E23").Value = 12.0850989817002
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 31.32788980611694
$ws.Range("H23").Value = 13.62848054262564
$ws.Range("M23").Value = 16.11703876236105
$ws.Range("N23").Value = 16.4016196407089
$ws.Range("B24").Value = 14.7026977905337
$ws.Range("C24").Value = 10.30078939586349
$ws.Range("E24").Value = 11.84845905833784
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 30.4444351167869
$ws.Range("H24").Value = 13.61422909575321
$ws.Range("M24").Value = 15.60000960025407
$ws.Range("N24").Value = 16.51887943334912
$ws.Range("B25").Value = 13.63845658638759
$ws.Range("C25").Value = 9.491851390287472
$ws.Range("E25").Value = 11.6027706136276
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 29.5186201445479
$ws.Range("H25").Value = 13.61398634533897
$ws.Range("M25").Value = 15.0376943418188
$ws.Range("N25").Value = 16.65270532493963
